$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

function Set-PlainCell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Row 2
Set-TextCell "D2" "302.54"
Set-TextCell "E2" "-0.69%"

# Row 3
Set-TextCell "D3" "36.77"
Set-TextCell "E3" "3.10%"

# Row 4
Set-TextCell "D4" "5.010"
Set-TextCell "E4" "-1.72%"

# Row 5
Set-TextCell "D5" "0.07699"
Set-TextCell "E5" "-1.30%"

# Row 6
Set-TextCell "D6" "2.088"
Set-TextCell "E6" "-7.83%"

# Row 7
Set-TextCell "D7" "7.999"
Set-TextCell "E7" "-1.40%"

# Row 8
Set-PlainCell "B8" "GateToken"
Set-PlainCell "C8" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell "D8" "4.050"
Set-TextCell "E8" "0.34%"

# Row 9
Set-PlainCell "B9" "MXToken"
Set-PlainCell "C9" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D9" "0.9195"
Set-TextCell "E9" "-0.93%"

# Row 10
Set-PlainCell "B10" "LiechtensteinCryptoassetsExchange"
Set-PlainCell "C10" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D10" "0.09793"
Set-TextCell "E10" "0.46%"

# Row 11
Set-PlainCell "B11" "WazirX"
Set-PlainCell "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D11" "0.1852"
Set-TextCell "E11" "1.33%"

# Row 12
Set-PlainCell "B12" "MandalaExchangeToken"
Set-PlainCell "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D12" "0.08640"
Set-TextCell "E12" "0.93%"

# Row 13
Set-PlainCell "B13" "BitrueCoin"
Set-PlainCell "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D13" "0.03616"
Set-TextCell "E13" "5.52%"

# Row 14
Set-PlainCell "B14" "BitMartToken"
Set-PlainCell "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D14" "0.09970"
Set-TextCell "E14" "0.21%"

# Row 15
Set-PlainCell "B15" "BitForexToken"
Set-PlainCell "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D15" "0.001476"
Set-TextCell "E15" "-0.38%"

# Row 16
Set-PlainCell "B16" "CoinExToken"
Set-PlainCell "C16" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell "D16" "0.04624"
Set-TextCell "E16" "-1.38%"

# Row 17
Set-PlainCell "B17" "TigerCash"
Set-PlainCell "C17" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D17" "0.005751"
Set-TextCell "E17" "0.27%"

# Row 18
Set-PlainCell "B18" "LEO"
Set-PlainCell "C18" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D18" "3.472"
Set-TextCell "E18" "-0.21%"

# Row 19
Set-PlainCell "B19" "BTSEToken"
Set-PlainCell "C19" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell "D19" "2.553"
Set-TextCell "E19" "17.18%"

# Row 20
Set-PlainCell "B20" "BitpandaEcosystemToken"
Set-PlainCell "C20" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell "D20" "0.3433"
Set-TextCell "E20" "-0.84%"

# Row 21
Set-PlainCell "B21" "ProBitToken"
Set-PlainCell "C21" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextCell "D21" "0.1332"
Set-TextCell "E21" "0.73%"

# Row 22
Set-PlainCell "B22" "MCDex"
Set-PlainCell "C22" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell "D22" "4.963"
Set-TextCell "E22" "9.20%"

# Row 23
Set-PlainCell "B23" "ZBToken"
Set-PlainCell "C23" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextCell "D23" "0.2215"
Set-TextCell "E23" "-1.17%"

# Row 24
Set-TextCell "D24" "0.005102"
Set-TextCell "E24" "12.33%"

# Row 25
Set-TextCell "E25" "-0.27%"

# Row 26
Set-TextCell "D26" "0.0001409"
Set-TextCell "E26" "8.22%"

# Row 39
Set-TextCell "D39" "0.01774"
Set-TextCell "E39" "0.22%"

# Row 40
Set-TextCell "D40" "0.04650"
Set-TextCell "E40" "-1.52%"

# Row 41
Set-TextCell "D41" "0.007697"
Set-TextCell "E41" "-2.37%"

# Row 42
Set-TextCell "D42" "0.1391"
Set-TextCell "E42" "-1.98%"

# Row 43
Set-TextCell "D43" "0.007991"
Set-TextCell "E43" "-0.22%"

# Row 44
Set-TextCell "D44" "0.002255"
Set-TextCell "E44" "-1.70%"

# Row 45
Set-TextCell "D45" "0.009890"
Set-TextCell "E45" "8.61%"

# Row 46
Set-TextCell "D46" "0.00006311"
Set-TextCell "E46" "1.80%"

# Row 47
Set-TextCell "E47" "0.48%"

# Row 48
Set-TextCell "D48" "0.0005831"
Set-TextCell "E48" "0.52%"

# Row 49
Set-TextCell "D49" "34.20"
Set-TextCell "E49" "518.81%"

# Row 50
Set-TextCell "D50" "0.002014"
Set-TextCell "E50" "-25.27%"

# Row 51
Set-TextCell "E51" "0.48%"
